$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 404.42856
$ws.Range("I4").Value = 76.59999999999999
$ws.Range("K4").Value = 76.59999999999999
$ws.Range("M4").Value = 37.40000000000001
$ws.Range("H9").Value = 3186.973
$ws.Range("I9").Value = 4352.6
$ws.Range("J9").Value = 758.5833
$ws.Range("K9").Value = 4352.6
$ws.Range("L9").Value = 758.5833
$ws.Range("M9").Value = -4183.6
$ws.Range("N9").Value = -1096.5833
$ws.Range("H58").Value = 4409.533
$ws.Range("J58").Value = 12696.4
$ws.Range("L58").Value = 38089.2
$ws.Range("N58").Value = -38389.2
$ws.Range("H116").Value = 6238.5386
$ws.Range("I116").Value = 6258.4165
$ws.Range("J116").Value = 6000
$ws.Range("K116").Value = 6258.4165
$ws.Range("L116").Value = 6000
$ws.Range("M116").Value = -2816.4165
$ws.Range("N116").Value = -12884

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 12699.125
$ws.Range("I4").Value = 70.8
$ws.Range("J4").Value = 33746.332
$ws.Range("K4").Value = 70.8
$ws.Range("L4").Value = 33746.332
$ws.Range("M4").Value = 45.2
$ws.Range("N4").Value = -33978.332
$ws.Range("H26").Value = 5178.4
$ws.Range("I26").Value = 3248
$ws.Range("K26").Value = 3248
$ws.Range("M26").Value = -2918
$ws.Range("H63").Value = 4395
$ws.Range("J63").Value = 4592.5
$ws.Range("L63").Value = 4592.5
$ws.Range("N63").Value = -5964.5
$ws.Range("H66").Value = 4395
$ws.Range("J66").Value = 4592.5
$ws.Range("L66").Value = 22962.5
$ws.Range("N66").Value = -29826.5
$ws.Range("H122").Value = 2397.9092
$ws.Range("I122").Value = 2121.862
$ws.Range("J122").Value = 4399.25
$ws.Range("K122").Value = 6365.586
$ws.Range("L122").Value = 13197.75
$ws.Range("M122").Value = -3915.586
$ws.Range("N122").Value = -18097.75
$ws.Range("H132").Value = 2780661.2
$ws.Range("I132").Value = 2943964.8
$ws.Range("K132").Value = 8831894.399999999
$ws.Range("M132").Value = -8829364.399999999

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3521
$ws.Range("I105").Value = 2900.7144
$ws.Range("K105").Value = 2900.7144
$ws.Range("M105").Value = -1153.7144

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1553806.2
$ws.Range("I16").Value = 1812440.6
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 1812440.6
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -1812153.6
$ws.Range("N16").Value = -2574
$ws.Range("H105").Value = 2001545.9
$ws.Range("I105").Value = 3334493.2
$ws.Range("J105").Value = 2125
$ws.Range("K105").Value = 3334493.2
$ws.Range("L105").Value = 2125
$ws.Range("M105").Value = -3332746.2
$ws.Range("N105").Value = -5619
$ws.Range("H113").Value = 1553806.2
$ws.Range("I113").Value = 1812440.6
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1812440.6
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = -1810270.6
$ws.Range("N113").Value = -6340
$ws.Range("H122").Value = 2972.72
$ws.Range("I122").Value = 2958.6191
$ws.Range("J122").Value = 3046.75
$ws.Range("K122").Value = 8875.8573
$ws.Range("L122").Value = 9140.25
$ws.Range("M122").Value = -6425.8573
$ws.Range("N122").Value = -14040.25

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10875.866
$ws.Range("I70").Value = 10840.125
$ws.Range("K70").Value = 10840.125
$ws.Range("M70").Value = -10570.125
$ws.Range("H73").Value = 10875.866
$ws.Range("I73").Value = 10840.125
$ws.Range("K73").Value = 10840.125
$ws.Range("M73").Value = -9904.125
$ws.Range("H80").Value = 3000
$ws.Range("J80").Value = 3000
$ws.Range("L80").Value = 3000
$ws.Range("N80").Value = -4996
$ws.Range("H83").Value = 3000
$ws.Range("J83").Value = 3000
$ws.Range("L83").Value = 15000
$ws.Range("N83").Value = -24984
$ws.Range("H99").Value = 12892.333
$ws.Range("I99").Value = 4570.9
$ws.Range("K99").Value = 4570.9
$ws.Range("M99").Value = -2324.9
$ws.Range("H132").Value = 5953901
$ws.Range("I132").Value = 6251346.5
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 18754039.5
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -18751509.5
$ws.Range("N132").Value = -20057
$ws.Range("H141").Value = 84764.39999999999
$ws.Range("J141").Value = 84764.39999999999
$ws.Range("L141").Value = 84764.39999999999
$ws.Range("N141").Value = -95124.39999999999

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
$ws.Range("H22").Value = 4277.6665
$ws.Range("I22").Value = 3833.3333
$ws.Range("J22").Value = 4499.8335
$ws.Range("K22").Value = 3833.3333
$ws.Range("L22").Value = 4499.8335
$ws.Range("M22").Value = -3538.3333
$ws.Range("N22").Value = -5089.8335
$ws.Range("H27").Value = 4277.6665
$ws.Range("I27").Value = 3833.3333
$ws.Range("J27").Value = 4499.8335
$ws.Range("K27").Value = 3833.3333
$ws.Range("L27").Value = 4499.8335
$ws.Range("M27").Value = -3726.3333
$ws.Range("N27").Value = -4713.8335
$ws.Range("H68").Value = 2977404.5
$ws.Range("I68").Value = 4762904
$ws.Range("J68").Value = 1571.3334
$ws.Range("K68").Value = 4762904
$ws.Range("L68").Value = 1571.3334
$ws.Range("M68").Value = -4762155
$ws.Range("N68").Value = -3069.3334
$ws.Range("H71").Value = 2977404.5
$ws.Range("I71").Value = 4762904
$ws.Range("J71").Value = 1571.3334
$ws.Range("K71").Value = 23814520
$ws.Range("L71").Value = 7856.666999999999
$ws.Range("M71").Value = -23810776
$ws.Range("N71").Value = -15344.667
$ws.Range("H93").Value = 500
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H101").Value = 25000
$ws.Range("J101").Value = 25000
$ws.Range("L101").Value = 25000
$ws.Range("N101").Value = -31490
$ws.Range("H132").Value = 19241210
$ws.Range("I132").Value = 20844458
$ws.Range("J132").Value = 2222
$ws.Range("K132").Value = 62533374
$ws.Range("L132").Value = 6666
$ws.Range("M132").Value = -62530844
$ws.Range("N132").Value = -11726

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 8333
$ws.Range("I5").Value = 9999
$ws.Range("J5").Value = 7500
$ws.Range("K5").Value = 9999
$ws.Range("L5").Value = 7500
$ws.Range("M5").Value = -9887
$ws.Range("N5").Value = -7724
$ws.Range("H100").Value = 1059
$ws.Range("I100").Value = 1052.1111
$ws.Range("J100").Value = 1090
$ws.Range("K100").Value = 2104.2222
$ws.Range("L100").Value = 2180
$ws.Range("M100").Value = -1563.2222
$ws.Range("N100").Value = -3262
$ws.Range("H122").Value = 8841.375
$ws.Range("I122").Value = 8841.375
$ws.Range("K122").Value = 26524.125
$ws.Range("M122").Value = -24074.125
$ws.Range("H136").Value = 25001570
$ws.Range("I136").Value = 35715132
$ws.Range("K136").Value = 107145396
$ws.Range("M136").Value = -107142846
$ws.Range("H140").Value = 72782.28999999999
$ws.Range("J140").Value = 72782.28999999999
$ws.Range("L140").Value = 72782.28999999999
$ws.Range("N140").Value = -83142.28999999999
$ws.Range("H141").Value = 77663
$ws.Range("I141").Value = 77000
$ws.Range("J141").Value = 77795.60000000001
$ws.Range("K141").Value = 77000
$ws.Range("L141").Value = 77795.60000000001
$ws.Range("M141").Value = -71820
$ws.Range("N141").Value = -88155.60000000001
